$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename (column headers -> snake_case) ---
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# --- Title-case Spanish connector words (de/del/la/las/los/el/y) in state/municipality names ---
$ws.Range('B7').Value = 'Pabellón De Arteaga'
$ws.Range('B8').Value = 'Rincón De Romos'
$ws.Range('B9').Value = 'San Francisco De Los Romo'
$ws.Range('B26').Value = 'Amatenango De La Frontera'
$ws.Range('B33').Value = 'Comitán De Domínguez'
$ws.Range('B41').Value = 'Mazapa De Madero'
$ws.Range('B50').Value = 'San Cristóbal De Las Casas'
$ws.Range('B83').Value = 'Guadalupe Y Calvo'
$ws.Range('B85').Value = 'Hidalgo Del Parral'
$ws.Range('B101').Value = 'San Francisco De Borja'
$ws.Range('B102').Value = 'San Francisco De Conchos'
$ws.Range('B103').Value = 'San Francisco Del Oro'
$ws.Range('B108').Value = 'Valle De Zaragoza'
$ws.Range('B130').Value = 'Villa De Álvarez'
$ws.Range('A132').Value = 'Ciudad De México'
$ws.Range('B136').Value = 'Cuajimalpa De Morelos'
$ws.Range('B150').Value = 'Coneto De Comonfort'
$ws.Range('B161').Value = 'Nombre De Dios'
$ws.Range('B164').Value = 'Pánuco De Coronado'
$ws.Range('B170').Value = 'San Luis Del Cordero'
$ws.Range('A179').Value = 'Estado De México'
$ws.Range('B179').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B182').Value = 'Almoloya De Alquisiras'
$ws.Range('B183').Value = 'Almoloya De Juárez'
$ws.Range('B187').Value = 'Atizapán De Zaragoza'
$ws.Range('B194').Value = 'Coacalco De Berriozábal'
$ws.Range('B199').Value = 'Ecatepec De Morelos'
$ws.Range('B203').Value = 'Ixtapan De La Sal'
$ws.Range('B214').Value = 'Naucalpan De Juárez'
$ws.Range('B221').Value = 'San Felipe Del Progreso'
$ws.Range('B222').Value = 'Soyaniquilpan De Juárez'
$ws.Range('B231').Value = 'Tenango Del Valle'
$ws.Range('B237').Value = 'Tlalnepantla De Baz'
$ws.Range('B243').Value = 'Valle De Bravo'
$ws.Range('B244').Value = 'Valle De Chalco Solidaridad'
$ws.Range('B253').Value = 'San Miguel De Allende'
$ws.Range('B254').Value = 'Apaseo El Alto'
$ws.Range('B255').Value = 'Apaseo El Grande'
$ws.Range('B259').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B263').Value = 'Jaral Del Progreso'
$ws.Range('B271').Value = 'Purísima Del Rincón'
$ws.Range('B275').Value = 'San Diego De La Unión'
$ws.Range('B277').Value = 'San Francisco Del Rincón'
$ws.Range('B279').Value = 'San Luis De La Paz'
$ws.Range('B280').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B282').Value = 'Silao De La Victoria'
$ws.Range('B286').Value = 'Valle De Santiago'
$ws.Range('B290').Value = 'Acapulco De Juárez'
$ws.Range('B293').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B294').Value = 'Alcozauca De Guerrero'
$ws.Range('B298').Value = 'Atenango Del Río'
$ws.Range('B299').Value = 'Atoyac De Álvarez'
$ws.Range('B300').Value = 'Ayutla De Los Libres'
$ws.Range('B303').Value = 'Buenavista De Cuéllar'
$ws.Range('B304').Value = 'Chilapa De Álvarez'
$ws.Range('B305').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B306').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B311').Value = 'Coyuca De Benítez'
$ws.Range('B312').Value = 'Coyuca De Catalán'
$ws.Range('B316').Value = 'Cuetzala Del Progreso'
$ws.Range('B317').Value = 'Cutzamala De Pinzón'
$ws.Range('B322').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B323').Value = 'Iguala De La Independencia'
$ws.Range('B325').Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range('B326').Value = 'Zihuatanejo De Azueta'
$ws.Range('B328').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B331').Value = 'Mártir De Cuilapan'
$ws.Range('B344').Value = 'Taxco De Alarcón'
$ws.Range('B346').Value = 'Técpan De Galeana'
$ws.Range('B348').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B350').Value = 'Tixtla De Guerrero'
$ws.Range('B353').Value = 'Tlalixtaquilla De Maldonado'
$ws.Range('B354').Value = 'Tlapa De Comonfort'
$ws.Range('B365').Value = 'Agua Blanca De Iturbide'
$ws.Range('B370').Value = 'Atotonilco El Grande'
$ws.Range('B374').Value = 'Cuautepec De Hinojosa'
$ws.Range('B378').Value = 'Huasca De Ocampo'
$ws.Range('B379').Value = 'Huejutla De Reyes'
$ws.Range('B381').Value = 'Jacala De Ledezma'
$ws.Range('B385').Value = 'Mineral Del Monte'
$ws.Range('B386').Value = 'Mixquiahuala De Juárez'
$ws.Range('B387').Value = 'Molango De Escamilla'
$ws.Range('B389').Value = 'Nopala De Villagrán'
$ws.Range('B390').Value = 'Pachuca De Soto'
$ws.Range('B393').Value = 'Progreso De Obregón'
$ws.Range('B395').Value = 'Santiago Tulantepec De Lugo Guerrero'
$ws.Range('B398').Value = 'Tenango De Doria'
$ws.Range('B400').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B405').Value = 'Tula De Allende'
$ws.Range('B406').Value = 'Tulancingo De Bravo'
$ws.Range('B411').Value = 'Acatlán De Juárez'
$ws.Range('B412').Value = 'Ahualulco De Mercado'
$ws.Range('B415').Value = 'Atemajac De Brizuela'
$ws.Range('B417').Value = 'Atotonilco El Alto'
$ws.Range('B419').Value = 'Autlán De Navarro'
$ws.Range('B430').Value = 'Concepción De Buenos Aires'
$ws.Range('B431').Value = 'Cuautitlán De García Barragán'
$ws.Range('B435').Value = 'Encarnación De Díaz'
$ws.Range('B441').Value = 'Huejuquilla El Alto'
$ws.Range('B442').Value = 'Ixtlahuacán Del Río'
$ws.Range('B446').Value = 'Jilotlán De Los Dolores'
$ws.Range('B451').Value = 'Lagos De Moreno'
$ws.Range('B462').Value = 'San Cristóbal De La Barranca'
$ws.Range('B463').Value = 'San Diego De Alejandría'
$ws.Range('B465').Value = 'San Juan De Los Lagos'
$ws.Range('B466').Value = 'San Juanito De Escobedo'
$ws.Range('B469').Value = 'San Miguel El Alto'
$ws.Range('B470').Value = 'San Sebastián Del Oeste'
$ws.Range('B471').Value = 'Santa María De Los Ángeles'
$ws.Range('B474').Value = 'Talpa De Allende'
$ws.Range('B475').Value = 'Tamazula De Gordiano'
$ws.Range('B481').Value = 'Teocuitatlán De Corona'
$ws.Range('B482').Value = 'Tepatitlán De Morelos'
$ws.Range('B485').Value = 'Tizapán El Alto'
$ws.Range('B486').Value = 'Tlajomulco De Zúñiga'
$ws.Range('B496').Value = 'Unión De San Antonio'
$ws.Range('B497').Value = 'Unión De Tula'
$ws.Range('B498').Value = 'Valle De Guadalupe'
$ws.Range('B499').Value = 'Valle De Juárez'
$ws.Range('B503').Value = 'Yahualica De González Gallo'
$ws.Range('B504').Value = 'Zacoalco De Torres'
$ws.Range('B507').Value = 'Zapotitlán De Vadillo'
$ws.Range('B508').Value = 'Zapotlán Del Rey'
$ws.Range('B509').Value = 'Zapotlán El Grande'
$ws.Range('B527').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B529').Value = 'Cojumatlán De Régules'
$ws.Range('B612').Value = 'Coatlán Del Río'
$ws.Range('B618').Value = 'Jonacatepec De Leandro Valle'
$ws.Range('B621').Value = 'Puente De Ixtla'
$ws.Range('B626').Value = 'Tetela Del Volcán'
$ws.Range('B628').Value = 'Tlaltizapán De Zapata'
$ws.Range('B636').Value = 'Amatlán De Cañas'
$ws.Range('B637').Value = 'Bahía De Banderas'
$ws.Range('B641').Value = 'Ixtlán Del Río'
$ws.Range('B645').Value = 'Santa María Del Oro'
$ws.Range('B657').Value = 'San Nicolás De Los Garza'
$ws.Range('B660').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B663').Value = 'Chalcatongo De Hidalgo'
$ws.Range('B665').Value = 'Constancia Del Rosario'
$ws.Range('B667').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B668').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B669').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B671').Value = 'Ixtlán De Juárez'
$ws.Range('B672').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B676').Value = 'Mártires De Tacubaya'
$ws.Range('B677').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B678').Value = 'Oaxaca De Juárez'
$ws.Range('B679').Value = 'Ocotlán De Morelos'
$ws.Range('B680').Value = 'Pinotepa De Don Luis'
$ws.Range('B681').Value = 'Putla Villa De Guerrero'
$ws.Range('B699').Value = 'San José Del Progreso'
$ws.Range('B703').Value = 'San Juan Bautista Lo De Soto'
$ws.Range('B736').Value = 'Santa Inés Del Monte'
$ws.Range('B737').Value = 'Santa Lucía Del Camino'
$ws.Range('B767').Value = 'Santo Domingo De Morelos'
$ws.Range('B773').Value = 'Teotitlán De Flores Magón'
$ws.Range('B774').Value = 'Tezoatlán De Segura Y Luna'
$ws.Range('B775').Value = 'Tlacolula De Matamoros'
$ws.Range('B776').Value = 'Villa De Tamazulápam Del Progreso'
$ws.Range('B777').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B779').Value = 'Villa Sola De Vega'
$ws.Range('B780').Value = 'Zimatlán De Álvarez'
$ws.Range('B808').Value = 'Cuayuca De Andrade'
$ws.Range('B816').Value = 'Izúcar De Matamoros'
$ws.Range('B822').Value = 'Los Reyes De Juárez'
$ws.Range('B829').Value = 'Palmar De Bravo'
$ws.Range('B840').Value = 'San Salvador El Verde'
$ws.Range('B843').Value = 'Tecali De Herrera'
$ws.Range('B848').Value = 'Tepatlaxco De Hidalgo'
$ws.Range('B851').Value = 'Tepexi De Rodríguez'
$ws.Range('B856').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B873').Value = 'Amealco De Bonfil'
$ws.Range('B874').Value = 'Cadereyta De Montes'
$ws.Range('B876').Value = 'Jalpan De Serra'
$ws.Range('B879').Value = 'Pinal De Amoles'
$ws.Range('B882').Value = 'San Juan Del Río'
$ws.Range('B895').Value = 'Ciudad Del Maíz'
$ws.Range('B902').Value = 'Mexquitic De Carmona'
$ws.Range('B908').Value = 'Santa María Del Río'
$ws.Range('B914').Value = 'Villa De Arriaga'
$ws.Range('B915').Value = 'Villa De Ramos'
$ws.Range('B953').Value = 'Nacozari De García'
$ws.Range('B970').Value = 'Jalpa De Méndez'
$ws.Range('B992').Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range('B996').Value = 'Papalotla De Xicohténcatl'
$ws.Range('B999').Value = 'Tepetitla De Lardizábal'
$ws.Range('B1001').Value = 'Tetla De La Solidaridad'
$ws.Range('B1017').Value = 'Amatlán De Los Reyes'
$ws.Range('B1023').Value = 'Boca Del Río'
$ws.Range('B1025').Value = 'Castillo De Teayo'
$ws.Range('B1026').Value = 'Cazones De Herrera'
$ws.Range('B1034').Value = 'Cosamaloapan De Carpio'
$ws.Range('B1042').Value = 'Hueyapan De Ocampo'
$ws.Range('B1043').Value = 'Ignacio De La Llave'
$ws.Range('B1045').Value = 'Ixhuatlán De Madero'
$ws.Range('B1051').Value = 'Juchique De Ferrer'
$ws.Range('B1056').Value = 'Martínez De La Torre'
$ws.Range('B1057').Value = 'Medellín De Bravo'
$ws.Range('B1067').Value = 'Paso De Ovejas'
$ws.Range('B1068').Value = 'Paso Del Macho'
$ws.Range('B1070').Value = 'Poza Rica De Hidalgo'
$ws.Range('B1076').Value = 'Sayula De Alemán'
$ws.Range('B1096').Value = 'Vega De Alatorre'
$ws.Range('B1116').Value = 'Cañitas De Felipe Pescador'
$ws.Range('B1117').Value = 'Concepción Del Oro'
$ws.Range('B1127').Value = 'Jiménez Del Teul'
$ws.Range('B1131').Value = 'Mezquital Del Oro'
$ws.Range('B1134').Value = 'Moyahua De Estrada'
$ws.Range('B1135').Value = 'Nochistlán De Mejía'
$ws.Range('B1136').Value = 'Noria De Ángeles'
$ws.Range('B1146').Value = 'Teúl De González Ortega'
$ws.Range('B1147').Value = 'Tlaltenango De Sánchez Román'
$ws.Range('B1149').Value = 'Villa De Cos'

# --- Floating point percentage value refresh (pct_matriculas recomputed with updated script) ---
$ws.Range('D11').Value = 0.009680429204571959
$ws.Range('D80').Value = 0.0009330534173081408
$ws.Range('D86').Value = 0.0009330534173081408
$ws.Range('D104').Value = 0.0009330534173081408
$ws.Range('D161').Value = 0.0009330534173081408
$ws.Range('D203').Value = 0.0009330534173081408
$ws.Range('D217').Value = 0.0009330534173081408
$ws.Range('D247').Value = 0.0009330534173081408
$ws.Range('D269').Value = 0.009447165850244929
$ws.Range('D309').Value = 0.0009330534173081408
$ws.Range('D332').Value = 0.0009330534173081408
$ws.Range('D333').Value = 0.0009330534173081408
$ws.Range('D373').Value = 0.0009330534173081408
$ws.Range('D459').Value = 0.0009330534173081408
$ws.Range('D475').Value = 0.0009330534173081408
$ws.Range('D502').Value = 0.0009330534173081408
$ws.Range('D507').Value = 0.0009330534173081408
$ws.Range('D523').Value = 0.0009330534173081408
$ws.Range('D546').Value = 0.0009330534173081408
$ws.Range('D621').Value = 0.0009330534173081408
$ws.Range('D628').Value = 0.0009330534173081408
$ws.Range('D637').Value = 0.0009330534173081408
$ws.Range('D755').Value = 0.0009330534173081408
$ws.Range('D762').Value = 0.0009330534173081408
$ws.Range('D772').Value = 0.0009330534173081408
$ws.Range('D880').Value = 0.0009330534173081408
$ws.Range('D934').Value = 0.0009330534173081408
$ws.Range('D982').Value = 0.0009330534173081408
$ws.Range('D1038').Value = 0.0009330534173081408
$ws.Range('D1087').Value = 0.0009330534173081408
$ws.Range('D1091').Value = 0.0009330534173081408
$ws.Range('D1147').Value = 0.0009330534173081408

# --- Remove trailing metadata/footer rows (1155:1159) ---
$ws.Rows('1155:1159').Delete()

